# Refresh cryptos list data (prices / 1h volume %) to the latest GitHub Actions run,
# including promoting Monero/Maker up a row and replacing Arweave with FLOKI at the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells so trailing zeros/precision are preserved
$textCells = @("D4", "D5", "D6", "D11", "D13", "D20", "D21", "D22", "D25", "D27", "D28", "D32", "D33", "D34", "D36", "D37", "D42", "D43", "D44", "D48", "D49", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range('D2').Value = '68.451.45'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '3.908.49'
$ws.Range('E3').Value = '  +3.65%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '602.13'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').Value = '164.45'
$ws.Range('E6').Value = '  -1.02%  '
$ws.Range('D7').Value = '3.909.90'
$ws.Range('E7').Value = '  +3.75%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  -2.43%  '
$ws.Range('E10').Value = '  -3.93%  '
$ws.Range('D11').Value = '6.38'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').Value = '36.79'
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').Value = '4.557.11'
$ws.Range('D16').Value = '3.877.20'
$ws.Range('E16').Value = '  +2.73%  '
$ws.Range('D17').Value = '68.636.58'
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  -1.32%  '
$ws.Range('D20').Value = '16.99'
$ws.Range('E20').Value = '  -3.70%  '
$ws.Range('D21').Value = '11.20'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').Value = '483.16'
$ws.Range('E22').Value = '  -2.12%  '
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('E24').Value = '  +11.77%  '
$ws.Range('D25').Value = '84.30'
$ws.Range('E25').Value = '  -0.64%  '
$ws.Range('E26').Value = '  -1.45%  '
$ws.Range('D27').Value = '11.96'
$ws.Range('E27').Value = '  -2.60%  '
$ws.Range('D28').Value = '10.08'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  -1.51%  '
$ws.Range('D31').Value = '4.056.05'
$ws.Range('E31').Value = '  +3.67%  '
$ws.Range('D32').Value = '7.83'
$ws.Range('E32').Value = '  -4.10%  '
$ws.Range('D33').Value = '2.36'
$ws.Range('E33').Value = '  -2.61%  '
$ws.Range('D34').Value = '31.86'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').Value = '3.849.99'
$ws.Range('E35').Value = '  +3.38%  '
$ws.Range('D36').Value = '0.106'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('D37').Value = '1.03'
$ws.Range('E37').Value = '  +2.07%  '
$ws.Range('E38').Value = '  +1.00%  '
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  -1.61%  '
$ws.Range('D42').Value = '0.316'
$ws.Range('E42').Value = '  -2.52%  '
$ws.Range('D43').Value = '431.84'
$ws.Range('E43').Value = '  +1.37%  '
$ws.Range('D44').Value = '48.46'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').Value = '26.36'
$ws.Range('E48').Value = '  +12.03%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '141.66'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.813.01'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').Value = '0.000264'
$ws.Range('E51').Value = '  +16.86%  '
